# "Double Bill Issue Fixed"
#
# This script reproduces, via Excel COM automation, the changes described
# by the target diff:
#   1. Stocks sheet: decrement on-hand Quantity (col D) for four items that
#      were sold in the new bills added below.
#   2. Bills sheet: normalize two stray "0.0" Bal_Amt text values to "0".
#   3. Bills sheet: append 7 new bill rows (157-163).
#
# Helper: the Bills sheet stores every column (including the numeric-looking
# Bill_Amt / Bal_Amt / Phone columns) as TEXT (shared strings), not as real
# numbers. Plain Excel automation auto-detects purely-numeric strings and
# stores them as numbers instead, so we force the "Text" number format
# before assigning the value and then restore the cell style to "Normal"
# (its original/default style) so no stray style index is left behind.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Stocks sheet - reduce Quantity (column D) for the items sold below
# ---------------------------------------------------------------------
$stocks = $wb.Worksheets.Item("Stocks")
$stocks.Cells.Item(2, 4).Value = 117.0   # Stock Name One : 139 -> 117 (-22)
$stocks.Cells.Item(4, 4).Value = 121.0   # NewStock_Img   : 132 -> 121 (-11)
$stocks.Cells.Item(5, 4).Value = 148.0   # Boost Half Kg  : 149 -> 148 (-1)
$stocks.Cells.Item(7, 4).Value = 197.0   # New Liz        : 199 -> 197 (-2)

# ---------------------------------------------------------------------
# 2. Bills sheet - fix the two stray "0.0" Bal_Amt clearance rows to "0"
# ---------------------------------------------------------------------
$bills = $wb.Worksheets.Item("Bills")
Set-TextValue $bills.Cells.Item(112, 5) "0"
Set-TextValue $bills.Cells.Item(156, 5) "0"

# ---------------------------------------------------------------------
# 3. Bills sheet - append the new bill rows 157-163
# ---------------------------------------------------------------------
# Row 157
$bills.Cells.Item(157, 1).Value = "15-Dec-2020 10:37"
$bills.Cells.Item(157, 2).Value = "ISM"
Set-TextValue $bills.Cells.Item(157, 3) "6587"
Set-TextValue $bills.Cells.Item(157, 4) "3410"
Set-TextValue $bills.Cells.Item(157, 5) "0"
$bills.Cells.Item(157, 6).Value = "XX1512153"
$bills.Cells.Item(157, 7).Value = "Stock Name One(22)"

# Row 158
$bills.Cells.Item(158, 1).Value = "15-Dec-2020 10:38"
$bills.Cells.Item(158, 2).Value = "ISM"
Set-TextValue $bills.Cells.Item(158, 3) "6587"
Set-TextValue $bills.Cells.Item(158, 4) "15"
Set-TextValue $bills.Cells.Item(158, 5) "0"
$bills.Cells.Item(158, 6).Value = "XX1512154"
$bills.Cells.Item(158, 7).Value = "New Liz(1)"

# Row 159
$bills.Cells.Item(159, 1).Value = "15-Dec-2020 10:38"
$bills.Cells.Item(159, 2).Value = "ISM"
Set-TextValue $bills.Cells.Item(159, 3) "6587"
Set-TextValue $bills.Cells.Item(159, 4) "1425"
Set-TextValue $bills.Cells.Item(159, 5) "0.0"
$bills.Cells.Item(159, 6).Value = "XX1512154"
$bills.Cells.Item(159, 7).Value = "New Liz(1),Bill Clearance 15Dec2020(1)"

# Row 160
$bills.Cells.Item(160, 1).Value = "15-Dec-2020 17:06"
$bills.Cells.Item(160, 2).Value = "Aakash"
Set-TextValue $bills.Cells.Item(160, 3) "364"
Set-TextValue $bills.Cells.Item(160, 4) "275"
Set-TextValue $bills.Cells.Item(160, 5) "0"
$bills.Cells.Item(160, 6).Value = "XX1512155"
$bills.Cells.Item(160, 7).Value = "NewStock_Img(11)"

# Row 161
$bills.Cells.Item(161, 1).Value = "15-Dec-2020 17:06"
$bills.Cells.Item(161, 2).Value = "Aakash"
Set-TextValue $bills.Cells.Item(161, 3) "364"
Set-TextValue $bills.Cells.Item(161, 4) "200"
Set-TextValue $bills.Cells.Item(161, 5) "0"
$bills.Cells.Item(161, 6).Value = "XX1512156"
$bills.Cells.Item(161, 7).Value = "Boost Half Kg(1)"

# Row 162
$bills.Cells.Item(162, 1).Value = "15-Dec-2020 17:07"
$bills.Cells.Item(162, 2).Value = "Aakash"
Set-TextValue $bills.Cells.Item(162, 3) "364"
Set-TextValue $bills.Cells.Item(162, 4) "375"
Set-TextValue $bills.Cells.Item(162, 5) "0"
$bills.Cells.Item(162, 6).Value = "XX1512156"
$bills.Cells.Item(162, 7).Value = "Boost Half Kg(1),Bill Clearance 15Dec2020(1)"

# Row 163
$bills.Cells.Item(163, 1).Value = "15-Dec-2020 17:12"
$bills.Cells.Item(163, 2).Value = "Aakash"
Set-TextValue $bills.Cells.Item(163, 3) "364"
Set-TextValue $bills.Cells.Item(163, 4) "340"
Set-TextValue $bills.Cells.Item(163, 5) "0.0"
$bills.Cells.Item(163, 6).Value = "XX1512157"
$bills.Cells.Item(163, 7).Value = "New Liz(1),Bill Clearance 15Dec2020(1)"
